$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.770.33'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.362.55'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '570.18'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.22%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '137.49'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.27%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.55%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '7.71'
$c.ClearFormats()
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('E11').Value = '  -3.80%  '
$ws.Range('D12').Value = '3.937.73'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('E13').Value = '  +0.64%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.91'
$c.ClearFormats()
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '3.363.80'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').Value = '60.901.92'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('E18').Value = '  -1.86%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.48'
$c.ClearFormats()
$ws.Range('E19').Value = '  -3.49%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '8.87'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.15%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '382.83'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '75.04'
$c.ClearFormats()
$ws.Range('E22').Value = '  +1.97%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.548'
$c.ClearFormats()
$ws.Range('E23').Value = '  -2.00%  '
$ws.Range('E24').Value = '  -0.14%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0000110'
$c.ClearFormats()
$ws.Range('E25').Value = '  -5.60%  '
$ws.Range('E26').Value = '  +6.04%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E27').Value = '  +0.34%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.11'
$c.ClearFormats()
$ws.Range('E28').Value = '  -3.82%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.88'
$c.ClearFormats()
$ws.Range('E30').Value = '  -1.93%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.34'
$c.ClearFormats()
$ws.Range('E32').Value = '  -6.49%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '22.95'
$c.ClearFormats()
$ws.Range('E33').Value = '  -2.65%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '6.81'
$c.ClearFormats()
$ws.Range('E34').Value = '  -2.10%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '166.97'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.48%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.92'
$c.ClearFormats()
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('D37').Value = '3.399.32'
$ws.Range('E37').Value = '  -0.42%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.43'
$c.ClearFormats()
$ws.Range('E38').Value = '  -3.17%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0753'
$c.ClearFormats()
$ws.Range('E39').Value = '  -2.43%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '25.52'
$c.ClearFormats()
$ws.Range('E40').Value = '  -8.47%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.770'
$c.ClearFormats()
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').Value = '2.436.78'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('E47').Value = '  -3.10%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '22.06'
$c.ClearFormats()
$ws.Range('E48').Value = '  -6.12%  '
$ws.Range('E49').Value = '  -5.33%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.97'
$c.ClearFormats()
$ws.Range('E50').Value = '  -4.68%  '
$ws.Range('E51').Value = '  -2.91%  '
